{"js": "// Replace each \"a+b=\" / \"a-b=\" arithmetic-problem cell with its new value.\n// Cells are visited in document (row-major) order; the Nth entry below\n// corresponds to the Nth table cell in the document.\nconst replacements = [[\"55+17=\", \"3+28=\"], [\"47-32=\", \"76-36=\"], [\"22-21=\", \"13+27=\"], [\"97-68=\", \"50+47=\"], [\"82-55=\", \"2+53=\"], [\"54+38=\", \"49-39=\"], [\"89+1=\", \"5+85=\"], [\"51-4=\", \"35-31=\"], [\"72+8=\", \"50-48=\"], [\"75-10=\", \"34+60=\"], [\"85+9=\", \"19-14=\"], [\"66+30=\", \"22-18=\"], [\"93-5=\", \"36-7=\"], [\"79-25=\", \"29-13=\"], [\"48-3=\", \"70-17=\"], [\"40+2=\", \"65-3=\"], [\"70-19=\", \"26-12=\"], [\"99-1=\", \"34+4=\"], [\"60-43=\", \"10+49=\"], [\"71-29=\", \"56-0=\"], [\"79-12=\", \"64+28=\"], [\"95-81=\", \"30+53=\"], [\"78-69=\", \"79-70=\"], [\"41+23=\", \"78-68=\"], [\"99-45=\", \"63+7=\"], [\"52-18=\", \"42-18=\"], [\"78-63=\", \"29+36=\"], [\"54-36=\", \"39-17=\"], [\"51-30=\", \"37+46=\"], [\"31+41=\", \"40+47=\"], [\"24+7=\", \"71-49=\"], [\"34+41=\", \"0+84=\"], [\"24+51=\", \"87-77=\"], [\"34+11=\", \"19+79=\"], [\"0+2=\", \"56+21=\"], [\"79-27=\", \"0+65=\"], [\"0+35=\", \"47-18=\"], [\"26-22=\", \"70-66=\"], [\"75+12=\", \"21+16=\"], [\"4+27=\", \"49+20=\"], [\"84-56=\", \"81-71=\"], [\"68+9=\", \"67-5=\"], [\"25+24=\", \"39-8=\"], [\"34+54=\", \"78-49=\"], [\"37+60=\", \"12+47=\"], [\"43+13=\", \"0+27=\"], [\"81-35=\", \"67-32=\"], [\"0+61=\", \"17+50=\"], [\"96-4=\", \"45+20=\"], [\"57-24=\", \"47+24=\"], [\"47+34=\", \"74+7=\"], [\"44+2=\", \"38-19=\"], [\"24+19=\", \"10+34=\"], [\"39+8=\", \"49+21=\"], [\"65+24=\", \"77+10=\"], [\"49-11=\", \"8-3=\"], [\"45+7=\", \"31+12=\"], [\"43+26=\", \"28-13=\"], [\"77+16=\", \"64-59=\"], [\"73-51=\", \"2+73=\"], [\"3+1=\", \"4+40=\"], [\"75+4=\", \"58+20=\"], [\"71-61=\", \"34+6=\"], [\"17+17=\", \"56-23=\"], [\"68+19=\", \"96-43=\"], [\"57-32=\", \"45-36=\"], [\"2+72=\", \"4+91=\"], [\"29+9=\", \"91-67=\"], [\"90+4=\", \"76-26=\"], [\"88-27=\", \"1+41=\"], [\"75-55=\", \"72-67=\"], [\"22-9=\", \"93-72=\"], [\"23+44=\", \"92-50=\"], [\"69-20=\", \"18+59=\"], [\"28-22=\", \"79-0=\"], [\"32+45=\", \"1+49=\"], [\"88-63=\", \"99-3=\"], [\"26+64=\", \"58+35=\"], [\"57-50=\", \"70-51=\"], [\"37+25=\", \"45-6=\"], [\"48-37=\", \"55+22=\"], [\"48-39=\", \"52+44=\"], [\"93-1=\", \"41+7=\"], [\"85-25=\", \"95-90=\"], [\"33+11=\", \"36-32=\"], [\"0+60=\", \"27-26=\"], [\"55-22=\", \"89-17=\"], [\"98-12=\", \"76-65=\"], [\"29-10=\", \"1+36=\"], [\"94-88=\", \"91-58=\"], [\"6+61=\", \"20+42=\"], [\"4+54=\", \"14+46=\"], [\"8+55=\", \"40-38=\"], [\"35+6=\", \"71-11=\"], [\"25+60=\", \"64-37=\"], [\"72-35=\", \"1+20=\"], [\"53+23=\", \"76-17=\"], [\"2+66=\", \"88-17=\"], [\"68-50=\", \"63+10=\"], [\"3+41=\", \"3+36=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Collect all cells in row-major order (the same order the cells appear in the OOXML).\nconst cells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cells.push(cell);\n  }\n}\n\nif (cells.length !== replacements.length) {\n  throw new Error(\n    \"Cell count (\" + cells.length + \") does not match replacement count (\" + replacements.length + \")\"\n  );\n}\n\n// Load each cell body's text up front so we can validate before mutating anything.\nfor (const cell of cells) {\n  cell.body.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cells.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const cell = cells[i];\n  const currentText = cell.body.text.replace(/\\r$/, \"\").trim();\n  if (currentText !== oldText) {\n    throw new Error(\n      \"Cell \" + i + \" text mismatch. Expected '\" + oldText + \"' but found '\" + currentText + \"'\"\n    );\n  }\n  // Use search() on the cell body so only the matching run's text is replaced,\n  // preserving the run's existing formatting (font, size, etc.).\n  const searchResults = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n  if (searchResults.items.length === 0) {\n    throw new Error(\"Could not find text '\" + oldText + \"' in cell \" + i);\n  }\n  searchResults.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each entry is the expected current text (\"old\") and its replacement (\"new\") for\n# one table cell, in row-major document order (same order the cells appear in the OOXML).\n$replacements = @(\n    @('55+17=', '3+28='),\n    @('47-32=', '76-36='),\n    @('22-21=', '13+27='),\n    @('97-68=', '50+47='),\n    @('82-55=', '2+53='),\n    @('54+38=', '49-39='),\n    @('89+1=', '5+85='),\n    @('51-4=', '35-31='),\n    @('72+8=', '50-48='),\n    @('75-10=', '34+60='),\n    @('85+9=', '19-14='),\n    @('66+30=', '22-18='),\n    @('93-5=', '36-7='),\n    @('79-25=', '29-13='),\n    @('48-3=', '70-17='),\n    @('40+2=', '65-3='),\n    @('70-19=', '26-12='),\n    @('99-1=', '34+4='),\n    @('60-43=', '10+49='),\n    @('71-29=', '56-0='),\n    @('79-12=', '64+28='),\n    @('95-81=', '30+53='),\n    @('78-69=', '79-70='),\n    @('41+23=', '78-68='),\n    @('99-45=', '63+7='),\n    @('52-18=', '42-18='),\n    @('78-63=', '29+36='),\n    @('54-36=', '39-17='),\n    @('51-30=', '37+46='),\n    @('31+41=', '40+47='),\n    @('24+7=', '71-49='),\n    @('34+41=', '0+84='),\n    @('24+51=', '87-77='),\n    @('34+11=', '19+79='),\n    @('0+2=', '56+21='),\n    @('79-27=', '0+65='),\n    @('0+35=', '47-18='),\n    @('26-22=', '70-66='),\n    @('75+12=', '21+16='),\n    @('4+27=', '49+20='),\n    @('84-56=', '81-71='),\n    @('68+9=', '67-5='),\n    @('25+24=', '39-8='),\n    @('34+54=', '78-49='),\n    @('37+60=', '12+47='),\n    @('43+13=', '0+27='),\n    @('81-35=', '67-32='),\n    @('0+61=', '17+50='),\n    @('96-4=', '45+20='),\n    @('57-24=', '47+24='),\n    @('47+34=', '74+7='),\n    @('44+2=', '38-19='),\n    @('24+19=', '10+34='),\n    @('39+8=', '49+21='),\n    @('65+24=', '77+10='),\n    @('49-11=', '8-3='),\n    @('45+7=', '31+12='),\n    @('43+26=', '28-13='),\n    @('77+16=', '64-59='),\n    @('73-51=', '2+73='),\n    @('3+1=', '4+40='),\n    @('75+4=', '58+20='),\n    @('71-61=', '34+6='),\n    @('17+17=', '56-23='),\n    @('68+19=', '96-43='),\n    @('57-32=', '45-36='),\n    @('2+72=', '4+91='),\n    @('29+9=', '91-67='),\n    @('90+4=', '76-26='),\n    @('88-27=', '1+41='),\n    @('75-55=', '72-67='),\n    @('22-9=', '93-72='),\n    @('23+44=', '92-50='),\n    @('69-20=', '18+59='),\n    @('28-22=', '79-0='),\n    @('32+45=', '1+49='),\n    @('88-63=', '99-3='),\n    @('26+64=', '58+35='),\n    @('57-50=', '70-51='),\n    @('37+25=', '45-6='),\n    @('48-37=', '55+22='),\n    @('48-39=', '52+44='),\n    @('93-1=', '41+7='),\n    @('85-25=', '95-90='),\n    @('33+11=', '36-32='),\n    @('0+60=', '27-26='),\n    @('55-22=', '89-17='),\n    @('98-12=', '76-65='),\n    @('29-10=', '1+36='),\n    @('94-88=', '91-58='),\n    @('6+61=', '20+42='),\n    @('4+54=', '14+46='),\n    @('8+55=', '40-38='),\n    @('35+6=', '71-11='),\n    @('25+60=', '64-37='),\n    @('72-35=', '1+20='),\n    @('53+23=', '76-17='),\n    @('2+66=', '88-17='),\n    @('68-50=', '63+10='),\n    @('3+41=', '3+36='),\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif (($rowCount * $colCount) -ne $replacements.Count) {\n    throw \"Cell count ($($rowCount * $colCount)) does not match replacement count ($($replacements.Count))\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $oldText = $replacements[$i][0]\n        $newText = $replacements[$i][1]\n\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        # Cell.Range includes the trailing end-of-cell marker; strip it for comparison.\n        $currentText = $cellRange.Text.TrimEnd([char]7, [char]13, [char]10)\n        if ($currentText -cne $oldText) {\n            throw \"Cell $i (row $r, col $c) text mismatch. Expected $oldText but found $currentText\"\n        }\n\n        $find = $cellRange.Find\n        $find.ClearFormatting()\n        $find.Text = $oldText\n        $find.Replacement.ClearFormatting()\n        $find.Replacement.Text = $newText\n        $find.Forward = $true\n        $find.Wrap = 0  # wdFindStop: stay within the cell range\n        $find.Format = $false\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.MatchWildcards = $false\n        $find.MatchSoundsLike = $false\n        $find.MatchAllWordForms = $false\n\n        $find.Execute(\n            [ref]$oldText, [ref]$true, [ref]$true, [ref]$false, [ref]$false,\n            [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$newText, [ref]2\n        ) | Out-Null\n\n        $i++\n    }\n}\n"}
